$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds free-form numeric-looking text (e.g. "1.0000",
# "0.9995", thousand-grouped "29.942.54"). Mark the cells being refreshed as
# Text first so Excel keeps the exact digit strings instead of normalising
# them into doubles (which would drop trailing zeros / alter precision).
$ws.Range("D2:D30,D32:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.942.54"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "1.892.82"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "0.8307"
$ws.Range("E5").Value = "  +8.38%  "
$ws.Range("D6").Value = "241.53"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("D7").Value = "1.000"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.3248"
$ws.Range("E8").Value = "  +6.65%  "
$ws.Range("D9").Value = "26.72"
$ws.Range("E9").Value = "  +5.40%  "
$ws.Range("D10").Value = "0.07027"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").Value = "0.08032"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "0.7480"
$ws.Range("D13").Value = "1.899.91"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "5.202"
$ws.Range("E14").Value = "  +0.66%  "
$ws.Range("D15").Value = "92.30"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "29.952.57"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "14.05"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "5.922"
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "243.65"
$ws.Range("E19").Value = "  -0.64%  "
$ws.Range("D20").Value = "0.000007759"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "2.151.75"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "1.0000"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "6.929"
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "0.1593"
$ws.Range("E25").Value = "  +24.42%  "
$ws.Range("D26").Value = "167.68"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "9.191"
$ws.Range("E27").Value = "  -0.73%  "
$ws.Range("D28").Value = "18.83"
$ws.Range("E28").Value = "  +0.88%  "
$ws.Range("D29").Value = "2.088"
$ws.Range("E29").Value = "  +2.63%  "
$ws.Range("D30").Value = "1.370"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "4.263"
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("D33").Value = "0.05636"
$ws.Range("E33").Value = "  +7.01%  "
$ws.Range("D34").Value = "4.072"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "1.277"
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").Value = "0.7323"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "0.01909"
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "2.779"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "0.4414"
$ws.Range("E40").Value = "  +0.38%  "
$ws.Range("D41").Value = "71.88"
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "5.940"
$ws.Range("E42").Value = "  -4.04%  "
$ws.Range("D43").Value = "0.8434"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "1.890"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("D46").Value = "101.44"
$ws.Range("E46").Value = "  +1.73%  "
$ws.Range("D47").Value = "7.602"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "9.749"
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "990.68"
$ws.Range("E49").Value = "  +9.01%  "
$ws.Range("D50").Value = "2.050.79"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("D51").Value = "36.05"
$ws.Range("E51").Value = "  -0.08%  "
